$wb = $excel.ActiveWorkbook

# --- Editor sheet: add two new test rows (Hii / NameError cases) ---
$ws2 = $wb.Worksheets.Item("Editor")

# Enter the new data in the same order it was typed by the author so the
# shared-string table is rebuilt with a matching order.
$ws2.Range("A4").Value = "Hii"
$ws2.Range("B4").Value = "NameError: name 'Hii' is not defined on line 1"
$ws2.Range("B5").Value = "NameError: empty not accepted"
$ws2.Range("A5").Value = " "

# Widen column B so the long error messages are fully visible.
$ws2.Columns.Item(2).ColumnWidth = 38.5

# --- Linkedlist sheet: move selection, no longer the active tab ---
$ws4 = $wb.Worksheets.Item("Linkedlist")
$ws4.Range("B2").Select() | Out-Null

# --- Login sheet: becomes the active tab with a new selection ---
$ws1 = $wb.Worksheets.Item("Login")
$ws1.Range("E14").Select() | Out-Null
